$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Cells.Item(2, 4) "64.305.44"
$ws.Range("E2").Value = "  +0.88%  "
Set-TextValue $ws.Cells.Item(3, 4) "2.764.73"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("E4").Value = "  +0.09%  "
Set-TextValue $ws.Cells.Item(5, 4) "577.32"
$ws.Range("E5").Value = "  +0.12%  "
Set-TextValue $ws.Cells.Item(6, 4) "160.94"
$ws.Range("E6").Value = "  +1.63%  "
Set-TextValue $ws.Cells.Item(7, 4) "0.999"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -1.05%  "
Set-TextValue $ws.Cells.Item(9, 4) "0.110"
$ws.Range("E9").Value = "  -0.79%  "
$ws.Range("E10").Value = "  +5.03%  "
Set-TextValue $ws.Cells.Item(11, 4) "5.86"
$ws.Range("E11").Value = "  +6.37%  "
Set-TextValue $ws.Cells.Item(12, 4) "0.387"
$ws.Range("E12").Value = "  -0.83%  "
Set-TextValue $ws.Cells.Item(13, 4) "3.257.20"
$ws.Range("E13").Value = "  +0.70%  "
Set-TextValue $ws.Cells.Item(14, 4) "27.45"
$ws.Range("E14").Value = "  +1.88%  "
Set-TextValue $ws.Cells.Item(15, 4) "63.956.55"
$ws.Range("E15").Value = "  +0.42%  "
Set-TextValue $ws.Cells.Item(16, 4) "0.0000152"
$ws.Range("E16").Value = "  -1.53%  "
Set-TextValue $ws.Cells.Item(17, 4) "2.773.54"
$ws.Range("E17").Value = "  +0.70%  "
Set-TextValue $ws.Cells.Item(18, 4) "12.21"
$ws.Range("E18").Value = "  -0.42%  "
Set-TextValue $ws.Cells.Item(19, 4) "4.86"
$ws.Range("E19").Value = "  -1.37%  "
Set-TextValue $ws.Cells.Item(20, 4) "358.53"
$ws.Range("E20").Value = "  -0.35%  "
Set-TextValue $ws.Cells.Item(21, 4) "6.69"
$ws.Range("E21").Value = "  -2.47%  "
$ws.Range("E22").Value = "  +0.79%  "
Set-TextValue $ws.Cells.Item(23, 4) "0.530"
$ws.Range("E23").Value = "  -6.09%  "
Set-TextValue $ws.Cells.Item(24, 4) "65.25"
$ws.Range("E24").Value = "  -1.61%  "
Set-TextValue $ws.Cells.Item(25, 4) "0.171"
$ws.Range("E25").Value = "  -0.28%  "
Set-TextValue $ws.Cells.Item(26, 4) "8.62"
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("E27").Value = "  -0.05%  "
Set-TextValue $ws.Cells.Item(28, 4) "0.0₃0924"
$ws.Range("E28").Value = "  -0.54%  "
Set-TextValue $ws.Cells.Item(29, 4) "7.38"
$ws.Range("E29").Value = "  +3.94%  "
Set-TextValue $ws.Cells.Item(30, 4) "1.99"
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("E31").Value = "  +10.64%  "
Set-TextValue $ws.Cells.Item(32, 4) "167.98"
$ws.Range("E32").Value = "  -0.73%  "
Set-TextValue $ws.Cells.Item(33, 4) "1.52"
$ws.Range("E33").Value = "  +4.39%  "
Set-TextValue $ws.Cells.Item(34, 4) "5.00"
$ws.Range("E34").Value = "  -0.09%  "
Set-TextValue $ws.Cells.Item(35, 4) "20.23"
$ws.Range("E35").Value = "  -1.22%  "
$ws.Range("E37").Value = "  +1.75%  "
Set-TextValue $ws.Cells.Item(38, 4) "1.01"
$ws.Range("E38").Value = "  +0.17%  "
Set-TextValue $ws.Cells.Item(39, 4) "350.71"
$ws.Range("E39").Value = "  +5.88%  "
$ws.Range("E40").Value = "  +5.10%  "
Set-TextValue $ws.Cells.Item(41, 4) "4.20"
$ws.Range("E41").Value = "  +0.18%  "
Set-TextValue $ws.Cells.Item(42, 4) "39.31"
$ws.Range("E42").Value = "  -0.78%  "
Set-TextValue $ws.Cells.Item(43, 4) "22.61"
$ws.Range("E43").Value = "  +3.49%  "
Set-TextValue $ws.Cells.Item(44, 4) "21.67"
$ws.Range("E44").Value = "  -1.28%  "
Set-TextValue $ws.Cells.Item(45, 4) "0.0596"
$ws.Range("E45").Value = "  +0.01%  "
Set-TextValue $ws.Cells.Item(46, 4) "136.91"
$ws.Range("E46").Value = "  +0.33%  "
Set-TextValue $ws.Cells.Item(47, 4) "0.633"
$ws.Range("E47").Value = "  -0.77%  "
Set-TextValue $ws.Cells.Item(48, 4) "0.0254"
$ws.Range("E48").Value = "  -1.43%  "
Set-TextValue $ws.Cells.Item(49, 4) "0.101"
$ws.Range("E49").Value = "  -0.37%  "
Set-TextValue $ws.Cells.Item(50, 4) "0.998"
$ws.Range("E50").Value = "  -0.15%  "
Set-TextValue $ws.Cells.Item(51, 4) "2.146.84"
$ws.Range("E51").Value = "  +1.65%  "
